$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrected/typo-fixed activity descriptions (rows 10-13, column C) ---
$ws.Range("C10").Value = "Actualización de teléfonos celulares. `nRelleno de toner en los contenedores de las impresoras"
$ws.Range("C11").Value = "Se realizo la migración de datos personales para reposición de equipo móvil.`nSe analizo la propuesta de realizar un sistema de lector de códigos de barras."
$ws.Range("C12").Value = "Se investigaron los programas a ocupar para el lector de códigos de barras, así como librerías, programas, y servidores."
$ws.Range("C13").Value = "Asignación y limpieza de equipos de cómputo.`nInvestigación de los lectores de códigos de barras que están en venta para poder llevar a cabo una propuesta de estos mismo."

# --- Realign the activity cells (C11:D13) from justified to centered text ---
$ws.Range("C11:D13").HorizontalAlignment = -4108

# --- Update the sheet view: zoom level and current selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 106
[void]$ws.Range("C13:D13").Select()
